$d = $word.ActiveDocument

# Fix "afficher" -> "affiché" (typo correction) in the first sentence.
$d.Content.Find.Execute("le mot afficher sur", $true, $false, $false, $false, $false,
                         $true, 1, $false, "le mot affiché sur", 2)

# Rework "le nombre de charactère" -> "les charactères inconnus"
$d.Content.Find.Execute("Les « _ » représente le nombre de charactère", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Les « _ » représente les charactères inconnus", 2)
